$wb = $excel.ActiveWorkbook

# --- Sheet "wheat": replace row 2 and append rows 3-6 ---
$wheat = $wb.Worksheets.Item("wheat")

$wheatData = @(
    @("FCSM", "Rajasthan", "DMSJ", "Bihar", "Wheat", 1),
    @("FCSM", "Rajasthan", "KSNG", "Odisha", "Wheat", 1),
    @("BGTN", "Punjab", "BBU", "Bihar", "Wheat", 1),
    @("KSA", "Punjab", "NNA", "Bihar", "Wheat", 1),
    @("JNL", "Punjab", "KSNG", "Odisha", "Wheat", 1)
)

$r = 2
foreach ($row in $wheatData) {
    $wheat.Cells.Item($r, 1).Value = $row[0]
    $wheat.Cells.Item($r, 2).Value = $row[1]
    $wheat.Cells.Item($r, 3).Value = $row[2]
    $wheat.Cells.Item($r, 4).Value = $row[3]
    $wheat.Cells.Item($r, 5).Value = $row[4]
    $wheat.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# --- Sheet "rra": remove row 2 data (only header remains) ---
$rra = $wb.Worksheets.Item("rra")
$rra.Rows.Item(2).Delete()
